$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated price (D) and 1h-volume-change (E) figures for the crypto list.
# D-column values that look like plain decimal numbers must be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# (losing formatting like trailing zeros, e.g. "20.00" -> 20).
$updates = @{
    2 = @{ D="41.308.53"; E="  -1.71%  " }
    3 = @{ D="2.181.15"; E="  -1.39%  " }
    4 = @{ E="  +0.07%  " }
    5 = @{ D="238.24"; E="  -0.86%  " }
    6 = @{ D="0.612"; E="  -2.01%  " }
    7 = @{ D="70.39"; E="  -3.54%  " }
    8 = @{ E="  +0.08%  " }
    9 = @{ E="  -4.05%  " }
    10 = @{ D="40.33"; E="  -5.89%  " }
    11 = @{ D="0.0929"; E="  -2.01%  " }
    12 = @{ D="54.32"; E="  -5.57%  " }
    14 = @{ E="  -4.02%  " }
    15 = @{ D="2.503.08"; E="  -1.41%  " }
    16 = @{ E="  -0.99%  " }
    17 = @{ D="0.803"; E="  -4.12%  " }
    18 = @{ D="2.175.27"; E="  -0.56%  " }
    19 = @{ D="41.131.73"; E="  -1.66%  " }
    20 = @{ E="  -5.84%  " }
    21 = @{ D="70.64"; E="  -2.77%  " }
    22 = @{ E="  -2.45%  " }
    23 = @{ D="9.79"; E="  -4.69%  " }
    24 = @{ D="226.51"; E="  -1.16%  " }
    25 = @{ E="  -5.86%  " }
    26 = @{ E="  +0.12%  " }
    27 = @{ D="10.92"; E="  -5.03%  " }
    28 = @{ E="  -0.74%  " }
    29 = @{ E="  -1.96%  " }
    30 = @{ E="  +0.54%  " }
    31 = @{ D="168.14"; E="  +0.52%  " }
    32 = @{ D="20.00"; E="  -2.64%  " }
    33 = @{ D="31.24"; E="  +7.60%  " }
    34 = @{ E="  -2.31%  " }
    35 = @{ D="5.16"; E="  -7.56%  " }
    36 = @{ E="  -2.97%  " }
    37 = @{ E="  -6.14%  " }
    38 = @{ E="  -3.16%  " }
    39 = @{ E="  -4.21%  " }
    40 = @{ E="  -1.24%  " }
    41 = @{ D="11.85"; E="  -8.43%  " }
    42 = @{ D="5.45"; E="  -2.77%  " }
    43 = @{ D="60.11"; E="  -7.39%  " }
    44 = @{ E="  -2.78%  " }
    45 = @{ E="  -2.51%  " }
    46 = @{ E="  -4.50%  " }
    47 = @{ D="98.58"; E="  -5.08%  " }
    48 = @{ E="  -1.57%  " }
    49 = @{ E="  -1.97%  " }
    50 = @{ E="  -7.11%  " }
    51 = @{ E="  -2.61%  " }
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    if ($vals.ContainsKey("D")) {
        $dCell = $ws.Range("D$row")
        $dVal = $vals["D"]
        if ($dVal -match '^-?[0-9]+(\.[0-9]+)?$') {
            # Looks numeric to Excel -- force Text format so it is stored
            # verbatim (matching the original inline-string cell) instead of
            # being coerced into a floating point number.
            $dCell.NumberFormat = "@"
        }
        $dCell.Value = $dVal
    }
    if ($vals.ContainsKey("E")) {
        $ws.Range("E$row").Value = $vals["E"]
    }
}
